$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Diebold-Mariano statistics (C) and p-values (D) per the commit
# "Correcion a Diebold Mariano y revision de Cap1"
$updates = @(
    @{ Row = 2;  C = -0.2744923139017411; D = 0.7862673926495143 },
    @{ Row = 3;  C = 0.7559947161944034;  D = 0.4576714854965753 },
    @{ Row = 4;  C = -0.4617498214073096; D = 0.6487941316509032 },
    @{ Row = 5;  C = -1.148978328098309;  D = 0.2629039614962512 },
    @{ Row = 6;  C = 0.9244875089613477;  D = 0.3652625026690521 },
    @{ Row = 7;  C = -0.2951176965389922; D = 0.7706711167417077 },
    @{ Row = 8;  C = -0.7860165257361048; D = 0.4402422524858056 },
    @{ Row = 9;  C = -0.9638842842096439; D = 0.3455801824238933 },
    @{ Row = 10; C = -1.649919001666828;  D = 0.1131658247870049 },
    @{ Row = 11; C = -0.4906242199881632; D = 0.6285502295590888 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}
